# "Bolão independencia - 130.xlsx" — trim the worksheet tab's name.
#
# The author moved the project files around (lotofacil/config.json,
# bolao-template.html, ...) and, while touching this workbook, shortened
# the single worksheet's tab name from "13 cotas(10 Jogos)" to
# "13 cotas". That is the only content-level change recorded for this
# workbook; the saved-window geometry and the cached
# "last known folder" hint that Excel stamps into the file are simply
# whatever the authoring machine's Excel session happened to be at save
# time and are not something a script drives deliberately.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab.
$ws.Name = "13 cotas"

# Best-effort extras matching the rest of the diff (harmless if the host
# doesn't persist window-geometry / last-path metadata):
# - the file was relocated under ...\Bolão\Bolao_Mega\lotofacil\
# - the saved window was moved/resized to xWindow=0 yWindow=0,
#   windowWidth=14610 windowHeight=15480
try { $wb.FullName = "C:\Users\allan\Desktop\Bolão\Bolao_Mega\lotofacil\Bolão independencia - 130.xlsx" } catch {}
try { $wb.Path = "C:\Users\allan\Desktop\Bolão\Bolao_Mega\lotofacil" } catch {}
try { $excel.ActiveWindow.Left = 0 } catch {}
try { $excel.ActiveWindow.Top = 0 } catch {}
try { $excel.ActiveWindow.Width = 14610 } catch {}
try { $excel.ActiveWindow.Height = 15480 } catch {}
